# Automated-task update: append the 15:00 reading row and correct the
# floating-point rounding of the previous (14:00) timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("14:00:24" reading): re-stamped serial date value (sub-second
# precision correction coming from the scheduled task's recompute).
$ws.Range("A10").Value = 45864.58361671296

# Give the new timestamp cell (A11) the same date style as the rest of
# column A by copying A10's formatting onto it before filling values in.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row 11: the 15:00 scheduled-task reading.
$ws.Range("A11").Value = 45864.62531725772
$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 21.44
$ws.Range("E11").Value = 67.06
$ws.Range("F11").Value = 425.39
$ws.Range("G11").Value = 14.55
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "15:00:27"
